$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values (columns B:AH) to 2 decimal places, matching Excel's ROUND() semantics
$cols = @(2..34)
foreach ($c in $cols) {
    $cell = $ws.Cells.Item(5, $c)
    $val = $cell.Value2
    if ($null -ne $val) {
        $d = [double]$val
        $sign = 1
        if ($d -lt 0) { $sign = -1 }
        $scaled = [Math]::Abs($d) * 100
        $rounded = [Math]::Floor($scaled + 0.5)
        $cell.Value2 = $sign * ($rounded / 100)
    }
}

# Delete row 6 entirely (shift cells up)
$ws.Rows.Item(6).Delete()
